# Edit script: Updated use cases and requirements
#
# 1) In the bullet "The system will allow the admin to add and remove
#    instructors.", change "remove" to "deactivate".
# 2) Add a new bullet at the end of the list:
#    "The system will provide access keys for use by students when a
#    class is unlocked."

$d = $word.ActiveDocument

# --- Change 1: remove -> deactivate -------------------------------------
$d.Content.Find.Execute("add and remove instructors", $true, $false, $false,
                         $false, $false, $true, 1, $false,
                         "add and deactivate instructors", 2)

# --- Change 2: add the new requirement bullet at the end of the list ----
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = "The system will provide access keys for use by students when a class is unlocked."
